# Update crypto prices / 1h-volume percentages (GitHub Actions refresh run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Price (column D), Volume(1h) (column E). $null means "leave as-is".
$updates = @(
    @{ Row = 2;  D = "67.273.44";  E = "  +0.45%  " },
    @{ Row = 3;  D = "3.482.91";   E = $null },
    @{ Row = 4;  D = $null;        E = "  -0.03%  " },
    @{ Row = 5;  D = $null;        E = "  -0.15%  " },
    @{ Row = 6;  D = "178.24";     E = "  +3.27%  " },
    @{ Row = 7;  D = $null;        E = "  +0.01%  " },
    @{ Row = 8;  D = "0.601";      E = "  +1.81%  " },
    @{ Row = 9;  D = $null;        E = "  -0.43%  " },
    @{ Row = 10; D = $null;        E = "  +4.59%  " },
    @{ Row = 11; D = $null;        E = "  -2.14%  " },
    @{ Row = 12; D = $null;        E = "  +0.62%  " },
    @{ Row = 13; D = $null;        E = "  -0.51%  " },
    @{ Row = 14; D = $null;        E = "  +9.14%  " },
    @{ Row = 16; D = $null;        E = "  +0.45%  " },
    @{ Row = 17; D = $null;        E = "  -1.21%  " },
    @{ Row = 18; D = $null;        E = "  -0.29%  " },
    @{ Row = 19; D = $null;        E = "  -0.78%  " },
    @{ Row = 20; D = $null;        E = "  -0.23%  " },
    @{ Row = 21; D = $null;        E = "  -1.90%  " },
    @{ Row = 22; D = $null;        E = "  +0.36%  " },
    @{ Row = 23; D = "73.85";      E = "  +0.48%  " },
    @{ Row = 24; D = $null;        E = "  +0.17%  " },
    @{ Row = 25; D = $null;        E = "  +0.24%  " },
    @{ Row = 26; D = $null;        E = "  +0.52%  " },
    @{ Row = 27; D = $null;        E = "  -0.62%  " },
    @{ Row = 28; D = $null;        E = "  +0.63%  " },
    @{ Row = 29; D = $null;        E = "  -4.10%  " },
    @{ Row = 30; D = $null;        E = "  +0.35%  " },
    @{ Row = 31; D = $null;        E = "  -1.25%  " },
    @{ Row = 32; D = $null;        E = "  -0.95%  " },
    @{ Row = 33; D = $null;        E = "  -0.43%  " },
    @{ Row = 34; D = $null;        E = "  -0.97%  " },
    @{ Row = 35; D = $null;        E = "  -0.43%  " },
    @{ Row = 36; D = $null;        E = "  +0.03%  " },
    @{ Row = 37; D = $null;        E = "  -1.81%  " },
    @{ Row = 38; D = $null;        E = "  +1.06%  " },
    @{ Row = 39; D = $null;        E = "  -1.44%  " },
    @{ Row = 40; D = $null;        E = "  +8.23%  " },
    @{ Row = 41; D = $null;        E = "  -2.52%  " },
    @{ Row = 42; D = $null;        E = "  -1.33%  " },
    @{ Row = 43; D = $null;        E = "  +0.02%  " },
    @{ Row = 44; D = $null;        E = "  -0.29%  " },
    @{ Row = 45; D = "27.02";      E = "  -0.39%  " },
    @{ Row = 46; D = $null;        E = "  -0.63%  " },
    @{ Row = 47; D = $null;        E = "  -2.46%  " },
    @{ Row = 48; D = $null;        E = "  -2.73%  " },
    @{ Row = 49; D = $null;        E = "  -0.93%  " },
    @{ Row = 50; D = $null;        E = "  -0.67%  " },
    @{ Row = 51; D = $null;        E = "  -2.51%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $origStyle = $cell.Style
        # Force the replacement text to stay a text value (matches the
        # existing inline-string cell type) instead of being auto-coerced
        # to a number by Excel's literal-input parsing.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = $origStyle
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
